# Edit: C1--C2-and-C3-PowerPoint.pptx
#
# 1) The table on slide 16 (the 3rd shape, a 2-column summary table) had its
#    table style switched from the custom "Table_0" style
#    ({3863D14D-A81D-4CB8-A6EE-C54A03DCC967}) to the built-in PowerPoint
#    table style {C162031E-F412-4022-A3F3-0D59D85D02AC}.
#
# 2) The presentation's design theme (used by the slide master / every
#    slide) was switched from the custom "Integral" palette to the
#    standard Office palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{C162031E-F412-4022-A3F3-0D59D85D02AC}")
}

# --- 2) Swap the design's colour scheme to the standard Office colours ----
# OOXML order: dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4,
#              accent5, accent6, hlink, folHlink
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $hex = $officeColors[$i]
    $r = ($hex -band 0xFF0000) -shr 16
    $g = ($hex -band 0x00FF00) -shr 8
    $b = ($hex -band 0x0000FF)
    $bgr = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $colorScheme.Colors($i + 1).RGB = $bgr
}
